$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82 (ALC)
$ws.Cells.Item(82, 8).Value = 793.3333
$ws.Cells.Item(82, 9).Value = 793.3333
$ws.Cells.Item(82, 11).Value = 2379.9999
$ws.Cells.Item(82, 13).Value = -1973.9999

# Row 85 (ALC)
$ws.Cells.Item(85, 8).Value = 793.3333
$ws.Cells.Item(85, 9).Value = 793.3333
$ws.Cells.Item(85, 11).Value = 2379.9999
$ws.Cells.Item(85, 13).Value = -975.9998999999998

# Row 88 (ALC)
$ws.Cells.Item(88, 8).Value = 5856.7144
$ws.Cells.Item(88, 9).Value = 4497
$ws.Cells.Item(88, 11).Value = 4497
$ws.Cells.Item(88, 13).Value = -4091

# Row 91 (ALC)
$ws.Cells.Item(91, 8).Value = 5856.7144
$ws.Cells.Item(91, 9).Value = 4497
$ws.Cells.Item(91, 11).Value = 4497
$ws.Cells.Item(91, 13).Value = -3093

# Row 106 (ALC)
$ws.Cells.Item(106, 8).Value = 8017.231
$ws.Cells.Item(106, 9).Value = 3261.6
$ws.Cells.Item(106, 10).Value = 14502.182
$ws.Cells.Item(106, 11).Value = 3261.6
$ws.Cells.Item(106, 12).Value = 14502.182
$ws.Cells.Item(106, 13).Value = -2630.6
$ws.Cells.Item(106, 14).Value = -15764.182

# Row 129 (ALC)
$ws.Cells.Item(129, 8).Value = 2484.2273
$ws.Cells.Item(129, 9).Value = 2199.7144
$ws.Cells.Item(129, 10).Value = 2617
$ws.Cells.Item(129, 11).Value = 6599.1432
$ws.Cells.Item(129, 12).Value = 7851
$ws.Cells.Item(129, 13).Value = -1599.1432
$ws.Cells.Item(129, 14).Value = -17851

# Row 132 (ALC)
$ws.Cells.Item(132, 8).Value = 5642.2
$ws.Cells.Item(132, 9).Value = 5802.75
$ws.Cells.Item(132, 11).Value = 17408.25
$ws.Cells.Item(132, 13).Value = -14878.25

# Row 135 (ALC)
$ws.Cells.Item(135, 8).Value = 822.3461
$ws.Cells.Item(135, 9).Value = 645.45
$ws.Cells.Item(135, 10).Value = 1412
$ws.Cells.Item(135, 11).Value = 5809.05
$ws.Cells.Item(135, 12).Value = 12708
$ws.Cells.Item(135, 13).Value = -3274.05
$ws.Cells.Item(135, 14).Value = -17778

# Row 138 (ALC)
$ws.Cells.Item(138, 8).Value = 3428.8604
$ws.Cells.Item(138, 10).Value = 3584
$ws.Cells.Item(138, 12).Value = 10752
$ws.Cells.Item(138, 14).Value = -21032

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Cells.Item(32, 8).Value = 2450.94
$ws.Cells.Item(32, 9).Value = 1558.9535
$ws.Cells.Item(32, 10).Value = 7930.2856
$ws.Cells.Item(32, 11).Value = 1558.9535
$ws.Cells.Item(32, 12).Value = 7930.2856
$ws.Cells.Item(32, 13).Value = -1271.9535
$ws.Cells.Item(32, 14).Value = -8504.285599999999

# Row 45 (ARM)
$ws.Cells.Item(45, 8).Value = 62504080
$ws.Cells.Item(45, 9).Value = 166667980
$ws.Cells.Item(45, 10).Value = 5737.3
$ws.Cells.Item(45, 11).Value = 166667980
$ws.Cells.Item(45, 12).Value = 5737.3
$ws.Cells.Item(45, 13).Value = -166667603
$ws.Cells.Item(45, 14).Value = -6491.3

# Row 61 (ARM)
$ws.Cells.Item(61, 8).Value = 3725.7368
$ws.Cells.Item(61, 9).Value = 2892.7334
$ws.Cells.Item(61, 11).Value = 2892.7334
$ws.Cells.Item(61, 13).Value = -2680.7334

# Row 122 (ARM)
$ws.Cells.Item(122, 8).Value = 1673.5625
$ws.Cells.Item(122, 9).Value = 1194.6428
$ws.Cells.Item(122, 11).Value = 3583.9284
$ws.Cells.Item(122, 13).Value = -1133.9284

# Row 132 (ARM)
$ws.Cells.Item(132, 8).Value = 8039.625
$ws.Cells.Item(132, 9).Value = 3436.3333
$ws.Cells.Item(132, 11).Value = 10308.9999
$ws.Cells.Item(132, 13).Value = -7778.999899999999

# Row 136 (ARM)
$ws.Cells.Item(136, 8).Value = 3725.7368
$ws.Cells.Item(136, 9).Value = 2892.7334
$ws.Cells.Item(136, 11).Value = 8678.200199999999
$ws.Cells.Item(136, 13).Value = -6128.200199999999

$ws = $wb.Worksheets.Item("BSM")
# Row 64 (BSM)
$ws.Cells.Item(64, 8).Value = 3222.111
$ws.Cells.Item(64, 10).Value = 3557
$ws.Cells.Item(64, 12).Value = 3557
$ws.Cells.Item(64, 14).Value = -4007

# Row 67 (BSM)
$ws.Cells.Item(67, 8).Value = 3222.111
$ws.Cells.Item(67, 10).Value = 3557
$ws.Cells.Item(67, 12).Value = 3557
$ws.Cells.Item(67, 14).Value = -5117

# Row 99 (BSM)
$ws.Cells.Item(99, 8).Value = 3080
$ws.Cells.Item(99, 10).Value = 4000
$ws.Cells.Item(99, 12).Value = 4000
$ws.Cells.Item(99, 14).Value = -6996

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (CRP)
$ws.Cells.Item(58, 8).Value = 6107.607
$ws.Cells.Item(58, 9).Value = 5257.857
$ws.Cells.Item(58, 11).Value = 5257.857
$ws.Cells.Item(58, 13).Value = -5054.857

# Row 94 (CRP)
$ws.Cells.Item(94, 8).Value = 2165.8125
$ws.Cells.Item(94, 9).Value = 1658.5
$ws.Cells.Item(94, 10).Value = 2334.9167
$ws.Cells.Item(94, 11).Value = 1658.5
$ws.Cells.Item(94, 12).Value = 2334.9167
$ws.Cells.Item(94, 13).Value = -1207.5
$ws.Cells.Item(94, 14).Value = -3236.9167

# Row 99 (CRP)
$ws.Cells.Item(99, 8).Value = 2337.5
$ws.Cells.Item(99, 9).Value = 2402.2
$ws.Cells.Item(99, 10).Value = 2014
$ws.Cells.Item(99, 11).Value = 2402.2
$ws.Cells.Item(99, 12).Value = 2014
$ws.Cells.Item(99, 13).Value = -904.1999999999998
$ws.Cells.Item(99, 14).Value = -5010

# Row 122 (CRP)
$ws.Cells.Item(122, 8).Value = 6662.25
$ws.Cells.Item(122, 9).Value = 2149.6
$ws.Cells.Item(122, 10).Value = 14183.333
$ws.Cells.Item(122, 11).Value = 6448.799999999999
$ws.Cells.Item(122, 12).Value = 42549.999
$ws.Cells.Item(122, 13).Value = -3998.799999999999
$ws.Cells.Item(122, 14).Value = -47449.999

# Row 126 (CRP)
$ws.Cells.Item(126, 8).Value = 2337.5
$ws.Cells.Item(126, 9).Value = 2402.2
$ws.Cells.Item(126, 10).Value = 2014
$ws.Cells.Item(126, 11).Value = 7206.599999999999
$ws.Cells.Item(126, 12).Value = 6042
$ws.Cells.Item(126, 13).Value = -4736.599999999999
$ws.Cells.Item(126, 14).Value = -10982

# Row 132 (CRP)
$ws.Cells.Item(132, 8).Value = 6240.269
$ws.Cells.Item(132, 9).Value = 6237.864
$ws.Cells.Item(132, 10).Value = 6253.5
$ws.Cells.Item(132, 11).Value = 18713.592
$ws.Cells.Item(132, 12).Value = 18760.5
$ws.Cells.Item(132, 13).Value = -16183.592
$ws.Cells.Item(132, 14).Value = -23820.5

# Row 134 (CRP)
$ws.Cells.Item(134, 8).Value = 3205.6667
$ws.Cells.Item(134, 9).Value = 2444.65
$ws.Cells.Item(134, 10).Value = 5380
$ws.Cells.Item(134, 11).Value = 7333.950000000001
$ws.Cells.Item(134, 12).Value = 16140
$ws.Cells.Item(134, 13).Value = -4798.950000000001
$ws.Cells.Item(134, 14).Value = -21210

# Row 136 (CRP)
$ws.Cells.Item(136, 8).Value = 6107.607
$ws.Cells.Item(136, 9).Value = 5257.857
$ws.Cells.Item(136, 11).Value = 15773.571
$ws.Cells.Item(136, 13).Value = -13223.571

$ws = $wb.Worksheets.Item("CUL")
# Row 61 (CUL)
$ws.Cells.Item(61, 8).Value = 1353.8889
$ws.Cells.Item(61, 9).Value = 169
$ws.Cells.Item(61, 10).Value = 1946.3334
$ws.Cells.Item(61, 11).Value = 507
$ws.Cells.Item(61, 12).Value = 5839.0002
$ws.Cells.Item(61, 13).Value = -292
$ws.Cells.Item(61, 14).Value = -6269.0002

# Row 62 (CUL)
$ws.Cells.Item(62, 8).Value = 19999
$ws.Cells.Item(62, 10).Value = 19999
$ws.Cells.Item(62, 12).Value = 59997
$ws.Cells.Item(62, 14).Value = -61369

# Row 65 (CUL)
$ws.Cells.Item(65, 8).Value = 19999
$ws.Cells.Item(65, 10).Value = 19999
$ws.Cells.Item(65, 12).Value = 179991
$ws.Cells.Item(65, 14).Value = -186855

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (GSM)
$ws.Cells.Item(102, 8).Value = 2015.5862
$ws.Cells.Item(102, 9).Value = 1435.1666
$ws.Cells.Item(102, 11).Value = 1435.1666
$ws.Cells.Item(102, 13).Value = 186.8334

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 152716.14
$ws.Cells.Item(132, 9).Value = 502499.5
$ws.Cells.Item(132, 10).Value = 12802.8
$ws.Cells.Item(132, 11).Value = 1507498.5
$ws.Cells.Item(132, 12).Value = 38408.39999999999
$ws.Cells.Item(132, 13).Value = -1504968.5
$ws.Cells.Item(132, 14).Value = -43468.39999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Cells.Item(7, 8).Value = 5036.387
$ws.Cells.Item(7, 9).Value = 2350.8635
$ws.Cells.Item(7, 11).Value = 2350.8635
$ws.Cells.Item(7, 13).Value = -2238.8635

# Row 16 (LTW)
$ws.Cells.Item(16, 8).Value = 3815.0833
$ws.Cells.Item(16, 9).Value = 3087
$ws.Cells.Item(16, 10).Value = 5999.3335
$ws.Cells.Item(16, 11).Value = 3087
$ws.Cells.Item(16, 12).Value = 5999.3335
$ws.Cells.Item(16, 13).Value = -2917
$ws.Cells.Item(16, 14).Value = -6339.3335

# Row 42 (LTW)
$ws.Cells.Item(42, 8).Value = 39012.5
$ws.Cells.Item(42, 9).Value = 39012.5
$ws.Cells.Item(42, 11).Value = 39012.5
$ws.Cells.Item(42, 13).Value = -38449.5

# Row 49 (LTW)
$ws.Cells.Item(49, 8).Value = 39012.5
$ws.Cells.Item(49, 9).Value = 39012.5
$ws.Cells.Item(49, 11).Value = 39012.5
$ws.Cells.Item(49, 13).Value = -38865.5

# Row 61 (LTW)
$ws.Cells.Item(61, 8).Value = 17005
$ws.Cells.Item(61, 10).Value = 17005
$ws.Cells.Item(61, 12).Value = 17005
$ws.Cells.Item(61, 14).Value = -17409

# Row 113 (LTW)
$ws.Cells.Item(113, 8).Value = 17005
$ws.Cells.Item(113, 10).Value = 17005
$ws.Cells.Item(113, 12).Value = 17005
$ws.Cells.Item(113, 14).Value = -21345

# Row 126 (LTW)
$ws.Cells.Item(126, 8).Value = 5036.387
$ws.Cells.Item(126, 9).Value = 2350.8635
$ws.Cells.Item(126, 11).Value = 7052.5905
$ws.Cells.Item(126, 13).Value = -4582.5905

$ws = $wb.Worksheets.Item("WVR")
# Row 44 (WVR)
$ws.Cells.Item(44, 8).Value = 23873.5
$ws.Cells.Item(44, 10).Value = 23873.5
$ws.Cells.Item(44, 12).Value = 23873.5
$ws.Cells.Item(44, 14).Value = -24981.5

# Row 104 (WVR)
$ws.Cells.Item(104, 8).Value = 12685
$ws.Cells.Item(104, 10).Value = 12685
$ws.Cells.Item(104, 12).Value = 12685
$ws.Cells.Item(104, 14).Value = -19673

# Row 126 (WVR)
$ws.Cells.Item(126, 8).Value = 1982.3077
$ws.Cells.Item(126, 9).Value = 1101.5
$ws.Cells.Item(126, 11).Value = 3304.5
$ws.Cells.Item(126, 13).Value = -834.5

# Row 132 (WVR)
$ws.Cells.Item(132, 8).Value = 4684.607
$ws.Cells.Item(132, 9).Value = 4423.5415
$ws.Cells.Item(132, 11).Value = 13270.6245
$ws.Cells.Item(132, 13).Value = -10740.6245
